$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.678.33'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '3.970.71'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'583.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.62%  '
$ws.Range("D6").Value = "'158.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.678"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.44%  '
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").Value = "'0.747"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").Value = "'53.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("E12").Value = '  -2.26%  '
$ws.Range("D13").Value = "'10.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("D14").Value = '4.604.99'
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").Value = '3.980.64'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("E16").Value = '  +7.25%  '
$ws.Range("D17").Value = "'13.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = "'20.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("D20").Value = '72.451.72'
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").Value = "'431.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = "'4.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.13%  '
$ws.Range("D23").Value = "'95.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.44%  '
$ws.Range("E24").Value = '  -3.65%  '
$ws.Range("D25").Value = "'14.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.04%  '
$ws.Range("D26").Value = "'4.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +22.86%  '
$ws.Range("D27").Value = "'11.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = "'10.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").Value = "'5.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").Value = "'36.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.28%  '
$ws.Range("D31").Value = "'7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.21%  '
$ws.Range("D32").Value = "'50.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.22%  '
$ws.Range("D33").Value = "'13.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").Value = "'678.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").Value = "'68.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.47%  '
$ws.Range("D37").Value = "'0.437"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '0.0₃0858'
$ws.Range("E38").Value = '  +3.13%  '
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -4.62%  '
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = "'10.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.99%  '
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = "'0.148"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("D47").Value = "'2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("E49").Value = '  +4.62%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("E51").Value = '  +6.63%  '
